{"js": "const replacements = [\n  [\"152\u00f72=76, 0\", \"319\u00f73=106, 1\"],\n  [\"342\u00f79=38, 0\", \"460\u00f74=115, 0\"],\n  [\"775\u00f72=387, 1\", \"200\u00f77=28, 4\"],\n  [\"928\u00f77=132, 4\", \"807\u00f74=201, 3\"],\n  [\"874\u00f78=109, 2\", \"182\u00f72=91, 0\"],\n  [\"713\u00f79=79, 2\", \"889\u00f75=177, 4\"],\n  [\"712\u00f73=237, 1\", \"547\u00f73=182, 1\"],\n  [\"542\u00f76=90, 2\", \"908\u00f76=151, 2\"],\n  [\"108\u00f72=54, 0\", \"585\u00f72=292, 1\"],\n  [\"896\u00f79=99, 5\", \"506\u00f75=101, 1\"],\n  [\"334\u00f73=111, 1\", \"497\u00f72=248, 1\"],\n  [\"587\u00f78=73, 3\", \"195\u00f78=24, 3\"],\n  [\"540\u00f77=77, 1\", \"991\u00f77=141, 4\"],\n  [\"721\u00f76=120, 1\", \"612\u00f73=204, 0\"],\n  [\"742\u00f74=185, 2\", \"989\u00f72=494, 1\"],\n  [\"180\u00f75=36, 0\", \"157\u00f75=31, 2\"],\n  [\"587\u00f74=146, 3\", \"899\u00f77=128, 3\"],\n  [\"342\u00f74=85, 2\", \"858\u00f77=122, 4\"],\n  [\"478\u00f79=53, 1\", \"370\u00f78=46, 2\"],\n  [\"248\u00f74=62, 0\", \"572\u00f76=95, 2\"],\n  [\"403\u00f77=57, 4\", \"789\u00f72=394, 1\"],\n  [\"766\u00f72=383, 0\", \"964\u00f77=137, 5\"],\n  [\"550\u00f79=61, 1\", \"613\u00f79=68, 1\"],\n  [\"298\u00f76=49, 4\", \"680\u00f77=97, 1\"],\n  [\"534\u00f77=76, 2\", \"999\u00f78=124, 7\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = \"152\u00f72=76, 0\"; New = \"319\u00f73=106, 1\"},\n    @{Old = \"342\u00f79=38, 0\"; New = \"460\u00f74=115, 0\"},\n    @{Old = \"775\u00f72=387, 1\"; New = \"200\u00f77=28, 4\"},\n    @{Old = \"928\u00f77=132, 4\"; New = \"807\u00f74=201, 3\"},\n    @{Old = \"874\u00f78=109, 2\"; New = \"182\u00f72=91, 0\"},\n    @{Old = \"713\u00f79=79, 2\"; New = \"889\u00f75=177, 4\"},\n    @{Old = \"712\u00f73=237, 1\"; New = \"547\u00f73=182, 1\"},\n    @{Old = \"542\u00f76=90, 2\"; New = \"908\u00f76=151, 2\"},\n    @{Old = \"108\u00f72=54, 0\"; New = \"585\u00f72=292, 1\"},\n    @{Old = \"896\u00f79=99, 5\"; New = \"506\u00f75=101, 1\"},\n    @{Old = \"334\u00f73=111, 1\"; New = \"497\u00f72=248, 1\"},\n    @{Old = \"587\u00f78=73, 3\"; New = \"195\u00f78=24, 3\"},\n    @{Old = \"540\u00f77=77, 1\"; New = \"991\u00f77=141, 4\"},\n    @{Old = \"721\u00f76=120, 1\"; New = \"612\u00f73=204, 0\"},\n    @{Old = \"742\u00f74=185, 2\"; New = \"989\u00f72=494, 1\"},\n    @{Old = \"180\u00f75=36, 0\"; New = \"157\u00f75=31, 2\"},\n    @{Old = \"587\u00f74=146, 3\"; New = \"899\u00f77=128, 3\"},\n    @{Old = \"342\u00f74=85, 2\"; New = \"858\u00f77=122, 4\"},\n    @{Old = \"478\u00f79=53, 1\"; New = \"370\u00f78=46, 2\"},\n    @{Old = \"248\u00f74=62, 0\"; New = \"572\u00f76=95, 2\"},\n    @{Old = \"403\u00f77=57, 4\"; New = \"789\u00f72=394, 1\"},\n    @{Old = \"766\u00f72=383, 0\"; New = \"964\u00f77=137, 5\"},\n    @{Old = \"550\u00f79=61, 1\"; New = \"613\u00f79=68, 1\"},\n    @{Old = \"298\u00f76=49, 4\"; New = \"680\u00f77=97, 1\"},\n    @{Old = \"534\u00f77=76, 2\"; New = \"999\u00f78=124, 7\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n"}
